# Auto-generated edit script: numeric market-price / profit-column corrections
# for the Seraph_Profits workbook (per-row Leve profit recalculation).
# Generated from the authoritative OOXML diff; cell refs are 1:1 with the
# sheet row each Leve/ingredient pair lives on.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 3728
$ws.Range("J88").Value = 4500
$ws.Range("L88").Value = 4500
$ws.Range("N88").Value = -5312
# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 3728
$ws.Range("J91").Value = 4500
$ws.Range("L91").Value = 4500
$ws.Range("N91").Value = -7308
# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 3067.0435
$ws.Range("I112").Value = 4200
$ws.Range("J112").Value = 2897.1
$ws.Range("K112").Value = 12600
$ws.Range("L112").Value = 8691.299999999999
$ws.Range("M112").Value = -11492
$ws.Range("N112").Value = -10907.3
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2001.7
$ws.Range("J137").Value = 1979.3334
$ws.Range("L137").Value = 5938.0002
$ws.Range("N137").Value = -11038.0002
# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 4373.8335
$ws.Range("I141").Value = 3898.8
$ws.Range("K141").Value = 11696.4
$ws.Range("M141").Value = -6516.400000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 1608.9474
$ws.Range("I45").Value = 1400.7273
$ws.Range("K45").Value = 1400.7273
$ws.Range("M45").Value = -1023.7273
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2203.577
$ws.Range("I132").Value = 2099.7827
$ws.Range("K132").Value = 6299.348100000001
$ws.Range("M132").Value = -3769.348100000001

$ws = $wb.Worksheets.Item("BSM")
# Row 112: Enlistment Highs / Deepgold Sword
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 2195
$ws.Range("I7").Value = 2321.8572
$ws.Range("K7").Value = 2321.8572
$ws.Range("M7").Value = -2208.8572
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 1251.3077
$ws.Range("I16").Value = 1225.4
$ws.Range("J16").Value = 1337.6666
$ws.Range("K16").Value = 1225.4
$ws.Range("L16").Value = 1337.6666
$ws.Range("M16").Value = -938.4000000000001
$ws.Range("N16").Value = -1911.6666
# Row 17: Say It with Spears / Feathered Harpoon
$ws.Range("H17").Value = 2800
$ws.Range("I17").Value = 2800
$ws.Range("K17").Value = 2800
$ws.Range("M17").Value = -2626
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 7039.615
$ws.Range("J31").Value = 9110
$ws.Range("L31").Value = 9110
$ws.Range("N31").Value = -9700
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 7039.615
$ws.Range("J34").Value = 9110
$ws.Range("L34").Value = 9110
$ws.Range("N34").Value = -9514
# Row 63: So You Think You Can Lance? / Mythrite Trident
$ws.Range("H63").Value = 70000
$ws.Range("J63").Value = 70000
$ws.Range("L63").Value = 70000
$ws.Range("N63").Value = -71372
# Row 66: Sticks and Stones (L) / Mythrite Trident
$ws.Range("H66").Value = 70000
$ws.Range("J66").Value = 70000
$ws.Range("L66").Value = 210000
$ws.Range("N66").Value = -216864
# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 1251.3077
$ws.Range("I113").Value = 1225.4
$ws.Range("J113").Value = 1337.6666
$ws.Range("K113").Value = 1225.4
$ws.Range("L113").Value = 1337.6666
$ws.Range("M113").Value = 944.5999999999999
$ws.Range("N113").Value = -5677.6666
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 5106.36
$ws.Range("I132").Value = 4155.0835
$ws.Range("K132").Value = 12465.2505
$ws.Range("M132").Value = -9935.250499999998
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1782.88
$ws.Range("I134").Value = 1332.8948
$ws.Range("K134").Value = 3998.6844
$ws.Range("M134").Value = -1463.6844

$ws = $wb.Worksheets.Item("CUL")
# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 20: Brothers in Arms / Brass Wristlets of Crafting
$ws.Range("H20").Value = 34888.43
$ws.Range("J20").Value = 34888.43
$ws.Range("L20").Value = 34888.43
$ws.Range("N20").Value = -35378.43
# Row 98: Cutting Deals / Durium Smallsword
$ws.Range("H98").Value = 13562.6
$ws.Range("J98").Value = 13562.6
$ws.Range("L98").Value = 13562.6
$ws.Range("N98").Value = -19552.6
# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 1878
$ws.Range("I107").Value = 1878
$ws.Range("K107").Value = 1878
$ws.Range("M107").Value = 42
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 32199.771
$ws.Range("J122").Value = 88024.336
$ws.Range("L122").Value = 264073.008
$ws.Range("N122").Value = -268973.008

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 3649
$ws.Range("I7").Value = 3694.818
$ws.Range("K7").Value = 3694.818
$ws.Range("M7").Value = -3582.818
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 1349.125
$ws.Range("J22").Value = 1399
$ws.Range("L22").Value = 1399
$ws.Range("N22").Value = -1989
# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 1349.125
$ws.Range("J27").Value = 1399
$ws.Range("L27").Value = 1399
$ws.Range("N27").Value = -1613
# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 2427.5
$ws.Range("I46").Value = 2375
$ws.Range("J46").Value = 2480
$ws.Range("K46").Value = 2375
$ws.Range("L46").Value = 2480
$ws.Range("M46").Value = -2187
$ws.Range("N46").Value = -2856
# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 990
$ws.Range("I55").Value = 990
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 990
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -817
$ws.Range("N55").ClearContents()
# Row 94: Fitting In / Gaganaskin Hat of Aiming
$ws.Range("H94").Value = 9999
$ws.Range("J94").Value = 9999
$ws.Range("L94").Value = 9999
$ws.Range("N94").Value = -11351
# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 3649
$ws.Range("I126").Value = 3694.818
$ws.Range("K126").Value = 11084.454
$ws.Range("M126").Value = -8614.454000000002

$ws = $wb.Worksheets.Item("WVR")
# Row 75: Storm upon Bald Mountain / Ramie Turban of Crafting
$ws.Range("H75").Value = 40118
$ws.Range("I75").Value = 40118
$ws.Range("K75").Value = 40118
$ws.Range("M75").Value = -39182
# Row 78: Abrupt Apprentices (L) / Ramie Turban of Crafting
$ws.Range("H78").Value = 40118
$ws.Range("I78").Value = 40118
$ws.Range("K78").Value = 120354
$ws.Range("M78").Value = -115674
# Row 80: Healing with Flair / Hallowed Ramie Gaskins of Healing
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 83: Pants Fit for Battle (L) / Hallowed Ramie Gaskins of Healing
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2522.0588
$ws.Range("J132").Value = 3050.25
$ws.Range("L132").Value = 9150.75
$ws.Range("N132").Value = -14210.75
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 4394.25
$ws.Range("I136").Value = 4394.25
$ws.Range("K136").Value = 13182.75
$ws.Range("M136").Value = -10632.75

Write-Host "Applied 155 value updates and 4 clears across 8 sheets."
